$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Retrain result: J1:K51 all become the constant 0.3 (was J1/K1 text "r"/"s",
# J2:J51 were 1, K2:K51 were already 0.3).
$ws.Range("J1:K51").Value = 0.3

# Update the view/selection state to match the saved workbook.
$ws.Range("K1:K51").Select()
